$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3-9 (keep header row 1 and data row 2)
$ws.Range("A3:H9").Delete() | Out-Null

# Update header row labels (category -> label/target naming)
$ws.Range("A1").Value = "labelA_stimuli"
$ws.Range("B1").Value = "labelB_stimuli"
$ws.Range("C1").Value = "targetA_stimuli"
$ws.Range("D1").Value = "targetB_stimuli"
$ws.Range("E1").Value = "labelA_image_stimuli"
$ws.Range("F1").Value = "labelB_image_stimuli"
$ws.Range("G1").Value = "targetA_image_stimuli"
$ws.Range("H1").Value = "targetB_image_stimuli"

# Update the selection/view state
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("A3:XFD10").Select() | Out-Null
$ws.Range("B3").Activate() | Out-Null
